# Insert a new worksheet "total_concentrations" right after
# "input_concentrations" (i.e. before "equilibrium_concentrations"),
# and populate it with the H / PO4 / Cu total-concentration values.

$wb = $excel.ActiveWorkbook

$insertAfter = $wb.Worksheets.Item("input_concentrations")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $insertAfter)
$newSheet.Name = "total_concentrations"

# Header row
$newSheet.Range("A1").Value = "H"
$newSheet.Range("B1").Value = "PO4"
$newSheet.Range("C1").Value = "Cu"

# Data rows
$newSheet.Range("A2").Value = 0.01
$newSheet.Range("B2").Value = 0.01
$newSheet.Range("C2").Value = 0.01

$newSheet.Range("A3").Value = 0.02
$newSheet.Range("B3").Value = 0.01
$newSheet.Range("C3").Value = 0.01

$newSheet.Range("A4").Value = 0.03
$newSheet.Range("B4").Value = 0.01
$newSheet.Range("C4").Value = 0.01

# Restore the originally-active sheet/tab so the workbook view is unchanged.
$wb.Worksheets.Item("input_stoich_coefficients").Activate()
